$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete appointment: remove the row for rajani@prac.to / slot 5 (original row 3),
# shifting the rows below it up by one.
$ws.Rows("3").Delete()

# Update what is now row 5 (previously row 6, uthmini@prac.to) with the new
# appointment details.
$ws.Range("B5").Value = "svbhadri1110@gmail.com"
$ws.Range("C5").Value = "18/04/2021"
$ws.Range("D5").Value = 1
